# Increase font sizes throughout the resume per the commit:
#   name header:        16pt -> 18pt
#   contact info:         9pt -> 10pt
#   section headers:     12pt -> 13pt
#   job titles:          11pt -> 12pt
#   body/bullets/dates:   9pt -> 10pt
#   overview paragraphs: 10pt -> 11pt
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    $rng = $d.Range($pr.Start, $pr.End)
    $size = $rng.Font.Size

    switch ($size) {
        16 { $rng.Font.Size = 18 }
        9  { $rng.Font.Size = 10 }
        12 { $rng.Font.Size = 13 }
        11 { $rng.Font.Size = 12 }
        10 { $rng.Font.Size = 11 }
        default { }
    }
}
